# ============================================================
# CustomWorldsInstalled.xlsx update
# Adds new "Source" column (E) with hyperlinks, two new worlds
# (Kingdom Hearts 1 / Re:Chain of Memories), refreshed version
# numbers + release/install dates, and widens columns A & E.
# ============================================================
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Date/number-format source cell (style already used by column C/D) ---
$dateFormatCell = $ws.Range("D2")

# --- Header row: add column E "Source" (bold/centered like B1:D1) ---
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value2 = "Source"

# --- Clear stale cells left over from rows whose World changed after re-sort ---
$ws.Range("C13").Clear()
$ws.Range("B14").Clear()
$ws.Range("B15").Clear()
$ws.Range("C16").Clear()

# --- Apply the date number-format to every Release/Install Date cell ---
# (PasteSpecial formats only, so existing style slot 5 is reused rather than
#  allocating a new numFmt; values are written separately below)
$dateFormatCell.Copy()
$ws.Range("C2").PasteSpecial(-4122)
$dateFormatCell.Copy()
$ws.Range("D2").PasteSpecial(-4122)
$dateFormatCell.Copy()
$ws.Range("C3").PasteSpecial(-4122)
$dateFormatCell.Copy()
$ws.Range("D3").PasteSpecial(-4122)
$dateFormatCell.Copy()
$ws.Range("D4").PasteSpecial(-4122)
$dateFormatCell.Copy()
$ws.Range("D5").PasteSpecial(-4122)
$dateFormatCell.Copy()
$ws.Range("C6").PasteSpecial(-4122)
$dateFormatCell.Copy()
$ws.Range("D6").PasteSpecial(-4122)
$dateFormatCell.Copy()
$ws.Range("C7").PasteSpecial(-4122)
$dateFormatCell.Copy()
$ws.Range("D7").PasteSpecial(-4122)
$dateFormatCell.Copy()
$ws.Range("D8").PasteSpecial(-4122)
$dateFormatCell.Copy()
$ws.Range("C9").PasteSpecial(-4122)
$dateFormatCell.Copy()
$ws.Range("D9").PasteSpecial(-4122)
$dateFormatCell.Copy()
$ws.Range("C10").PasteSpecial(-4122)
$dateFormatCell.Copy()
$ws.Range("D10").PasteSpecial(-4122)
$dateFormatCell.Copy()
$ws.Range("C11").PasteSpecial(-4122)
$dateFormatCell.Copy()
$ws.Range("D11").PasteSpecial(-4122)
$dateFormatCell.Copy()
$ws.Range("C12").PasteSpecial(-4122)
$dateFormatCell.Copy()
$ws.Range("D12").PasteSpecial(-4122)
$dateFormatCell.Copy()
$ws.Range("D13").PasteSpecial(-4122)
$dateFormatCell.Copy()
$ws.Range("D14").PasteSpecial(-4122)
$dateFormatCell.Copy()
$ws.Range("C15").PasteSpecial(-4122)
$dateFormatCell.Copy()
$ws.Range("D15").PasteSpecial(-4122)
$dateFormatCell.Copy()
$ws.Range("D16").PasteSpecial(-4122)
$dateFormatCell.Copy()
$ws.Range("C17").PasteSpecial(-4122)
$dateFormatCell.Copy()
$ws.Range("D17").PasteSpecial(-4122)
$dateFormatCell.Copy()
$ws.Range("C18").PasteSpecial(-4122)
$dateFormatCell.Copy()
$ws.Range("D18").PasteSpecial(-4122)
$dateFormatCell.Copy()
$ws.Range("C19").PasteSpecial(-4122)
$dateFormatCell.Copy()
$ws.Range("D19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- World rows, already written in the final A-Z sorted order ---
# Row 2: A Hat in Time
$ws.Range("A2").Value2 = "A Hat in Time"
$ws.Range("B2").Value2 = "1.5.3"
$ws.Range("C2").Value2 = 45430
$ws.Range("D2").Value2 = 45435
$ws.Range("E2").Value2 = "https://github.com/CookieCat45/Archipelago-ahit/releases"

# Row 3: A Robot Named Fight!
$ws.Range("A3").Value2 = "A Robot Named Fight!"
$ws.Range("C3").Value2 = 45274
$ws.Range("D3").Value2 = 45422
$ws.Range("E3").Value2 = "https://discord.com/channels/731205301247803413/1169389087371841708"

# Row 4: Celeste
$ws.Range("A4").Value2 = "Celeste"
$ws.Range("D4").Value2 = 45410
$ws.Range("E4").Value2 = "https://github.com/doshyw/CelesteArchipelago/releases"

# Row 5: Chrono Trigger Jets of Time
$ws.Range("A5").Value2 = "Chrono Trigger Jets of Time"
$ws.Range("D5").Value2 = 45422
$ws.Range("E5").Value2 = "https://wiki.ctjot.com/doku.php?id=multiworld"

# Row 6: CrossCode
$ws.Range("A6").Value2 = "CrossCode"
$ws.Range("B6").Value2 = "0.4.4"
$ws.Range("C6").Value2 = 45428
$ws.Range("D6").Value2 = 45435
$ws.Range("E6").Value2 = "https://github.com/CodeTriangle/CCMultiworldRandomizer/releases"

# Row 7: Final Fantasy 12 Open World
$ws.Range("A7").Value2 = "Final Fantasy 12 Open World"
$ws.Range("B7").Value2 = "0.3.4"
$ws.Range("C7").Value2 = 45421
$ws.Range("D7").Value2 = 45435
$ws.Range("E7").Value2 = "https://github.com/Bartz24/Archipelago/releases"

# Row 8: Final Fantasy 5 Career Day
$ws.Range("A8").Value2 = "Final Fantasy 5 Career Day"
$ws.Range("B8").Value2 = "0.62"
$ws.Range("D8").Value2 = 45422
$ws.Range("E8").Value2 = "https://github.com/cleartonic/arch_ffvcd/releases"

# Row 9: Final Fantasy 6 Worlds Collide
$ws.Range("A9").Value2 = "Final Fantasy 6 Worlds Collide"
$ws.Range("C9").Value2 = 45428
$ws.Range("D9").Value2 = 45435
$ws.Range("E9").Value2 = "https://discord.com/channels/731205301247803413/1022545979146252288"

# Row 10: Inscryption
$ws.Range("A10").Value2 = "Inscryption"
$ws.Range("B10").Value2 = "0.2.0b2"
$ws.Range("C10").Value2 = 45265
$ws.Range("D10").Value2 = 45422
$ws.Range("E10").Value2 = "https://github.com/DrBibop/Archipelago_Inscryption/releases"

# Row 11: Kingdom Hearts 1
$ws.Range("A11").Value2 = "Kingdom Hearts 1"
$ws.Range("B11").Value2 = "v1.2.5"
$ws.Range("C11").Value2 = 45434
$ws.Range("D11").Value2 = 45435
$ws.Range("E11").Value2 = "https://github.com/gaithernOrg/KH1FM-AP/releases"

# Row 12: Kingdom Hearts Re:Chain of Memories
$ws.Range("A12").Value2 = "Kingdom Hearts Re:Chain of Memories"
$ws.Range("B12").Value2 = "v1.2.1"
$ws.Range("C12").Value2 = 45369
$ws.Range("D12").Value2 = 45435
$ws.Range("E12").Value2 = "https://github.com/gaithernOrg/ArchipelagoKHRECOM/releases"

# Row 13: Minit
$ws.Range("A13").Value2 = "Minit"
$ws.Range("B13").Value2 = "0.6.3"
$ws.Range("D13").Value2 = 45422
$ws.Range("E13").Value2 = "https://github.com/qwint/APMinit/releases"

# Row 14: Oracle of Seasons
$ws.Range("A14").Value2 = "Oracle of Seasons"
$ws.Range("D14").Value2 = 45410
$ws.Range("E14").Value2 = "https://github.com/Dinopony/ArchipelagoOoS/releases"

# Row 15: Pharcryption
$ws.Range("A15").Value2 = "Pharcryption"
$ws.Range("C15").Value2 = 45317
$ws.Range("D15").Value2 = 45422
$ws.Range("E15").Value2 = "https://discord.com/channels/731205301247803413/1092478908022136876"

# Row 16: Pseudoregalia
$ws.Range("A16").Value2 = "Pseudoregalia"
$ws.Range("B16").Value2 = "0.7.2"
$ws.Range("D16").Value2 = 45422
$ws.Range("E16").Value2 = "https://github.com/pseudoregalia-modding/pseudoregalia-archipelago/tags"

# Row 17: Shahrazad
$ws.Range("A17").Value2 = "Shahrazad"
$ws.Range("B17").Value2 = "0.1.0"
$ws.Range("C17").Value2 = 45414
$ws.Range("D17").Value2 = 45422
$ws.Range("E17").Value2 = "https://github.com/qwint/ap-shahrazad/releases"

# Row 18: Super Metroid Map Rando
$ws.Range("A18").Value2 = "Super Metroid Map Rando"
$ws.Range("B18").Value2 = "v111"
$ws.Range("C18").Value2 = 45429
$ws.Range("D18").Value2 = 45435
$ws.Range("E18").Value2 = "https://discord.com/channels/731205301247803413/1156395911874875473"

# Row 19: Wargroove 2
$ws.Range("A19").Value2 = "Wargroove 2"
$ws.Range("C19").Value2 = 45428
$ws.Range("D19").Value2 = 45435
$ws.Range("E19").Value2 = "https://discord.com/channels/731205301247803413/1159482310652076082"

# --- Refresh the persisted AutoSort state (data is already alphabetised) ---
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2:A19"))
$ws.Sort.SetRange($ws.Range("A2:D19"))
$ws.Sort.Header = -4163
$ws.Sort.Apply()

# --- Turn column E entries into real hyperlinks (adds Hyperlink style + rels) ---
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/CookieCat45/Archipelago-ahit/releases")
$ws.Hyperlinks.Add($ws.Range("E3"), "https://discord.com/channels/731205301247803413/1169389087371841708")
$ws.Hyperlinks.Add($ws.Range("E4"), "https://github.com/doshyw/CelesteArchipelago/releases")
$ws.Hyperlinks.Add($ws.Range("E5"), "https://wiki.ctjot.com/doku.php?id=multiworld")
$ws.Hyperlinks.Add($ws.Range("E6"), "https://github.com/CodeTriangle/CCMultiworldRandomizer/releases")
$ws.Hyperlinks.Add($ws.Range("E7"), "https://github.com/Bartz24/Archipelago/releases")
$ws.Hyperlinks.Add($ws.Range("E8"), "https://github.com/cleartonic/arch_ffvcd/releases")
$ws.Hyperlinks.Add($ws.Range("E9"), "https://discord.com/channels/731205301247803413/1022545979146252288")
$ws.Hyperlinks.Add($ws.Range("E10"), "https://github.com/DrBibop/Archipelago_Inscryption/releases")
$ws.Hyperlinks.Add($ws.Range("E11"), "https://github.com/gaithernOrg/KH1FM-AP/releases")
$ws.Hyperlinks.Add($ws.Range("E12"), "https://github.com/gaithernOrg/ArchipelagoKHRECOM/releases")
$ws.Hyperlinks.Add($ws.Range("E13"), "https://github.com/qwint/APMinit/releases")
$ws.Hyperlinks.Add($ws.Range("E14"), "https://github.com/Dinopony/ArchipelagoOoS/releases")
$ws.Hyperlinks.Add($ws.Range("E15"), "https://discord.com/channels/731205301247803413/1092478908022136876")
$ws.Hyperlinks.Add($ws.Range("E16"), "https://github.com/pseudoregalia-modding/pseudoregalia-archipelago/tags")
$ws.Hyperlinks.Add($ws.Range("E17"), "https://github.com/qwint/ap-shahrazad/releases")
$ws.Hyperlinks.Add($ws.Range("E18"), "https://discord.com/channels/731205301247803413/1156395911874875473")
$ws.Hyperlinks.Add($ws.Range("E19"), "https://discord.com/channels/731205301247803413/1159482310652076082")

# --- Column widths: A grew to fit the longest world name, E sized for URLs ---
$ws.Columns.Item(1).ColumnWidth = 34.592447916666664
$ws.Columns.Item(5).ColumnWidth = 74.16666666666667

# --- Selection / view bookkeeping ---
$ws.Range("E4").Select()
